$p = $ppt.ActivePresentation
$old = "21/10/2020"
$new = "26/10/2020"

function Update-ShapeDate {
    param($shapes)
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $sh = $shapes.Item($k)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq $old) {
                    $sh.TextFrame.TextRange.Text = $new
                }
            }
        }
    }
}

# Total number of slide layouts across every design/master is exposed
# globally through CustomLayouts.Item(), even though each design's
# CustomLayouts.Count only reports the slides that belong to it - so sum
# the per-design counts to get the real upper bound.
$totalLayouts = 0
for ($i = 1; $i -le $p.Designs.Count; $i++) {
    $design = $p.Designs.Item($i)
    $totalLayouts = $totalLayouts + $design.SlideMaster.CustomLayouts.Count
}

$firstMaster = $p.Designs.Item(1).SlideMaster
for ($j = 1; $j -le $totalLayouts; $j++) {
    $layout = $firstMaster.CustomLayouts.Item($j)
    Update-ShapeDate $layout.Shapes
}

# Update the date placeholder cached text on every slide master too.
for ($i = 1; $i -le $p.Designs.Count; $i++) {
    $design = $p.Designs.Item($i)
    $master = $design.SlideMaster
    Update-ShapeDate $master.Shapes
}
